# Weekly update: insert a new week's worth of Chirimoya price data
# (Especial / Primera / Segunda) at the top of the "recent" block,
# pushing the older rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows starting at row 312 (existing rows 312+ shift down by 3).
$ws.Range("A312:A314").EntireRow.Insert()

# --- Row 312: Especial ---
$ws.Range("A312").Value = 8
$ws.Range("B312").Value = "Terminal La Palmera de La Serena"
$ws.Range("C312").Value = "Coquimbo"
$ws.Range("D312").Value = 45223
$ws.Range("E312").Value = 4
$ws.Range("F312").Value = "Fruta"
$ws.Range("G312").Value = 100107
$ws.Range("H312").Value = "Otros"
$ws.Range("I312").Value = 100107002
$ws.Range("J312").Value = "Chirimoya"
$ws.Range("K312").Value = "Cultivar IV Región"
$ws.Range("L312").Value = "Especial"
$ws.Range("M312").Value = 300
$ws.Range("N312").Value = 20000
$ws.Range("O312").Value = 21000
$ws.Range("P312").Value = 20500
$ws.Range("Q312").Value = "`$/bandeja 10 kilos"
$ws.Range("R312").Value = "Provincia de Limarí"
$ws.Range("S312").Value = 2050
$ws.Range("T312").Value = 10

# --- Row 313: Primera ---
$ws.Range("A313").Value = 8
$ws.Range("B313").Value = "Terminal La Palmera de La Serena"
$ws.Range("C313").Value = "Coquimbo"
$ws.Range("D313").Value = 45223
$ws.Range("E313").Value = 4
$ws.Range("F313").Value = "Fruta"
$ws.Range("G313").Value = 100107
$ws.Range("H313").Value = "Otros"
$ws.Range("I313").Value = 100107002
$ws.Range("J313").Value = "Chirimoya"
$ws.Range("K313").Value = "Cultivar IV Región"
$ws.Range("L313").Value = "Primera"
$ws.Range("M313").Value = 240
$ws.Range("N313").Value = 17000
$ws.Range("O313").Value = 18000
$ws.Range("P313").Value = 17500
$ws.Range("Q313").Value = "`$/bandeja 10 kilos"
$ws.Range("R313").Value = "Provincia de Limarí"
$ws.Range("S313").Value = 1750
$ws.Range("T313").Value = 10

# --- Row 314: Segunda ---
$ws.Range("A314").Value = 8
$ws.Range("B314").Value = "Terminal La Palmera de La Serena"
$ws.Range("C314").Value = "Coquimbo"
$ws.Range("D314").Value = 45223
$ws.Range("E314").Value = 4
$ws.Range("F314").Value = "Fruta"
$ws.Range("G314").Value = 100107
$ws.Range("H314").Value = "Otros"
$ws.Range("I314").Value = 100107002
$ws.Range("J314").Value = "Chirimoya"
$ws.Range("K314").Value = "Cultivar IV Región"
$ws.Range("L314").Value = "Segunda"
$ws.Range("M314").Value = 200
$ws.Range("N314").Value = 14000
$ws.Range("O314").Value = 15000
$ws.Range("P314").Value = 14500
$ws.Range("Q314").Value = "`$/bandeja 10 kilos"
$ws.Range("R314").Value = "Provincia de Limarí"
$ws.Range("S314").Value = 1450
$ws.Range("T314").Value = 10
